$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Deskcount")

# Deskcount value corrections
$ws.Range("C12").Value = 79
$ws.Range("C43").Value = 32
$ws.Range("C44").Value = 561

# "Include in Occupancy Calculation" flips from Yes -> No
$ws.Range("F23").Value = "No"
$ws.Range("F46").Value = "No"
$ws.Range("F47").Value = "No"
$ws.Range("F48").Value = "No"

# Reflect the author's last on-screen scroll position / selection
$ws.Activate()
$excel.ActiveWindow.ScrollRow = 31
$excel.ActiveWindow.ScrollColumn = 1
$ws.Range("C42").Select() | Out-Null
